$words = @(
    'führen',
    'klingen',
    'scheinen',
    'trauen',
    'geben',
    'seufzen',
    'heben',
    'spüren',
    'schwingen',
    'schwören',
    'freuen',
    'räumen',
    'ehren',
    'platzen',
    'arten',
    'ändern',
    'schmecken',
    'warnen',
    'dringen',
    'bellen',
    'feiern',
    'wachsen',
    'flüchten',
    'schlucken',
    'wundern',
    'kehren',
    'sinken',
    'sterben',
    'heulen',
    'rasen',
    'boxen',
    'schwächen',
    'graben',
    'erben',
    'werden',
    'backen',
    'bluten',
    'sprengen',
    'fangen',
    'zünden',
    'sichern',
    'münzen',
    'stehlen',
    'werfen',
    'enden',
    'hauen',
    'stammen',
    'spinnen',
    'fällen',
    'jagen',
    'schenken',
    'betteln',
    'machen',
    'lesen',
    'reizen',
    'drehen',
    'brauchen',
    'irren',
    'zeigen',
    'gelten',
    'siegen',
    'planen',
    'biegen',
    'quälen',
    'hören',
    'helfen',
    'pflanzen',
    'decken',
    'wirken',
    'sorgen',
    'tropfen',
    'ärgern',
    'suchen',
    'töten',
    'grüßen',
    'wüten',
    'runden',
    'trennen',
    'achten',
    'altern',
    'pfeifen',
    'saufen',
    'schulden',
    'malen',
    'zielen',
    'scheitern',
    'mögen',
    'fahren',
    'zögern',
    'folgen',
    'liefern',
    'bauen',
    'schrecken',
    'treiben',
    'fließen',
    'rufen',
    'jubeln',
    'kosten',
    'lügen',
    'greifen',
    'streichen',
    'heilen',
    'loben',
    'dienen',
    'liegen',
    'knarren',
    'stecken',
    'kichern',
    'klettern',
    'tollen',
    'äußern',
    'wehtun',
    'gründen',
    'formen',
    'fallen',
    'lockern',
    'mauern',
    'spielen',
    'fischen',
    'schreiten',
    'filmen',
    'flehen',
    'bitten',
    'weichen',
    'bergen',
    'wenden',
    'kümmern',
    'sperren'
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $words[$i]
}
